$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy formatting from the last existing data row (25) down into the new row 26
# so the percent-format (C:W) / general-format (X:CG) split carries over.
$ws.Range("A25:CG25").Copy()
$ws.Range("A26").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Geography identifiers (new shared strings get created in this order so the
# full MSA name lands before the short name, matching the source workbook).
$ws.Range("B26").Value = "Detroit-Warren-Dearborn, MI"
$ws.Range("A26").Value = "Detroit"

# Quarterly inflation-rate series for the new Detroit MSA row.
$ws.Range("C26").Value = [double]"8.6891904761904795E-2"
$ws.Range("D26").Value = [double]"8.6371129707112998E-2"
$ws.Range("E26").Value = [double]"8.6690476190476207E-2"
$ws.Range("F26").Value = [double]"8.5960947259565698E-2"
$ws.Range("G26").Value = [double]"8.5721197593706597E-2"
$ws.Range("H26").Value = [double]"8.5898243520782402E-2"
$ws.Range("I26").Value = [double]"8.3731374382022494E-2"
$ws.Range("J26").Value = [double]"8.1858987662015506E-2"
$ws.Range("K26").Value = [double]"7.9192903162629794E-2"
$ws.Range("L26").Value = [double]"8.2126537588185694E-2"
$ws.Range("M26").Value = [double]"8.1587570275976698E-2"
$ws.Range("N26").Value = [double]"8.1097931034482804E-2"
$ws.Range("O26").Value = [double]"8.2752000000000006E-2"
$ws.Range("P26").Value = [double]"7.5851884751772999E-2"
$ws.Range("Q26").Value = [double]"7.55788494319663E-2"
$ws.Range("R26").Value = [double]"7.3813330662476404E-2"
$ws.Range("S26").Value = [double]"7.2618262944524906E-2"
$ws.Range("T26").Value = [double]"7.1114571914322605E-2"
$ws.Range("U26").Value = [double]"6.3463380923450802E-2"
$ws.Range("V26").Value = [double]"6.4546459226739294E-2"
$ws.Range("W26").Value = [double]"6.4111119958863302E-2"
$ws.Range("X26").Value = [double]"6.6044938344163098E-2"
$ws.Range("Y26").Value = [double]"7.1835993481562693E-2"
$ws.Range("Z26").Value = [double]"7.3208113375948E-2"
$ws.Range("AA26").Value = [double]"7.16589857753792E-2"
$ws.Range("AB26").Value = [double]"7.04586162836119E-2"
$ws.Range("AC26").Value = [double]"6.9022702702119104E-2"
$ws.Range("AD26").Value = [double]"6.9839038081795396E-2"
$ws.Range("AE26").Value = [double]"7.0523367109019394E-2"
$ws.Range("AF26").Value = [double]"8.2441144239219402E-2"
$ws.Range("AG26").Value = [double]"8.8201009824501894E-2"
$ws.Range("AH26").Value = [double]"8.7507053072625698E-2"
$ws.Range("AI26").Value = [double]"9.3173649289099505E-2"
$ws.Range("AJ26").Value = [double]"8.3526306620209101E-2"
$ws.Range("AK26").Value = [double]"8.4341598360655703E-2"
$ws.Range("AO26").Value = [double]"9.7564737031700294E-2"
$ws.Range("AP26").Value = [double]"9.67728858768407E-2"
$ws.Range("AQ26").Value = [double]"9.4790932657926097E-2"
$ws.Range("AR26").Value = [double]"9.3302271321614602E-2"
$ws.Range("AS26").Value = [double]"8.7079580796296296E-2"
$ws.Range("AT26").Value = [double]"7.8796806427469099E-2"
$ws.Range("AU26").Value = [double]"7.9067534542961102E-2"
$ws.Range("AV26").Value = [double]"8.0851605434038396E-2"
$ws.Range("AW26").Value = [double]"8.0098265332576093E-2"
$ws.Range("AX26").Value = [double]"8.0677228562061207E-2"
$ws.Range("AY26").Value = [double]"7.6309386043858701E-2"
$ws.Range("AZ26").Value = [double]"7.4428639481987399E-2"
$ws.Range("BA26").Value = [double]"7.4964094324016795E-2"
$ws.Range("BB26").Value = [double]"7.6374748748685606E-2"
$ws.Range("BC26").Value = [double]"7.8222792133929506E-2"
$ws.Range("BD26").Value = [double]"7.8412909051586197E-2"
$ws.Range("BE26").Value = [double]"7.7991645803927298E-2"
$ws.Range("BF26").Value = [double]"7.6313219436533297E-2"
$ws.Range("BG26").Value = [double]"7.5694098656561501E-2"
$ws.Range("BH26").Value = [double]"7.1021712280196506E-2"
$ws.Range("BI26").Value = [double]"7.0740044890639495E-2"
$ws.Range("BJ26").Value = [double]"6.9253260109880296E-2"
$ws.Range("BK26").Value = [double]"6.99945580345714E-2"
$ws.Range("BL26").Value = [double]"7.1317152145560794E-2"
$ws.Range("BM26").Value = [double]"7.0721272084129996E-2"
$ws.Range("BN26").Value = [double]"6.9993685425742594E-2"
$ws.Range("BO26").Value = [double]"6.7966677749999996E-2"
$ws.Range("BP26").Value = [double]"6.7112034695529801E-2"
$ws.Range("BQ26").Value = [double]"6.62331317652462E-2"
$ws.Range("BR26").Value = [double]"6.5757119355109597E-2"
$ws.Range("BS26").Value = [double]"6.6101983576761197E-2"
$ws.Range("BT26").Value = [double]"6.4935880315140707E-2"
$ws.Range("BU26").Value = [double]"6.4642532874853706E-2"
$ws.Range("BV26").Value = [double]"6.2988313062208201E-2"
$ws.Range("BW26").Value = [double]"6.3302016067390804E-2"
$ws.Range("BX26").Value = [double]"6.1417481676923097E-2"
$ws.Range("BY26").Value = [double]"6.1157200167224299E-2"
$ws.Range("BZ26").Value = [double]"6.0647290356573999E-2"
$ws.Range("CA26").Value = [double]"5.9340264784722201E-2"
$ws.Range("CB26").Value = [double]"5.9711754204394701E-2"
$ws.Range("CC26").Value = [double]"5.7166004197399103E-2"
$ws.Range("CD26").Value = [double]"5.6232553773217897E-2"
$ws.Range("CE26").Value = [double]"5.4767993074784298E-2"
$ws.Range("CF26").Value = [double]"5.4040305682230498E-2"
$ws.Range("CG26").Value = [double]"5.4783708443934501E-2"

# Leave AL26:AN26 and CI26 blank (no published BLS CPI yet for those quarters)
# and CH26 blank, matching the source data.

$ws.Range("A26").Select()
